# Adding new keys to localization file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("command0", "cd - change directory", "cd - zmień katalog"),
    @("command1", "ls - list files", "ls - wyświetl listę plików"),
    @("command2", "ssh - connect to device via ssh protocol", "ssh - połącz się z urządzeniem za pomocą protokołu ssh"),
    @("command3", "scp - copy files from remote machines via ssh", "scp - skopiuj pliki ze zdalnych maszyn za pomocą ssh"),
    @("command4", "mkdir - create directory", "mkdir - utwórz katalog")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 2).WrapText = $true

    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 3).WrapText = $true

    $row = $row + 1
}

$ws.Range("B13").Select() | Out-Null
